# Add a "Row Number" column (F) to Sheet1 that numbers each data row for
# graphing purposes (commit: "added row number for graphing").
#
# Column F gets:
#   F1            -> header "Row Number" (new shared string)
#   F2            -> standalone formula  =ROW(F2) - 1        (value 1)
#   F3:F66        -> shared formula group =ROW(F3) - 1        (values 2..65)
#   F67:F102      -> shared formula group =ROW(F67) - 1        (values 66..101)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header
$ws.Range("F1").Value = "Row Number"

# First data row gets its own (non-shared) formula.
$ws.Range("F2").Formula = "=ROW(F2) - 1"

# Remaining data rows, filled in two batches (mirrors how the source
# workbook ended up with two shared-formula groups: F3:F66 and F67:F102).
$ws.Range("F3:F66").Formula = "=ROW(F3) - 1"
$ws.Range("F67:F102").Formula = "=ROW(F67) - 1"

# Widen the new column to fit its header/content.
$ws.Columns.Item(6).ColumnWidth = 12.451822916666666

# Leave the new column selected, like the author did after filling it in.
$ws.Range("F2:F102").Select()
